# Generate Report for Handback
# Fills in the "handed back" status/target/handback columns for both the
# zh-cn and de-de localization sheets, mirrors the new status onto the
# Overview tab, and widens the columns that now hold the longer text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$mdFile1 = "28c8dcfe-f72d-4a1e-8572-60038800e9d0.md"
$mdFile2 = "ca40efca-d1b0-4c5f-af3c-95fc671a0aee.md"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/10e063f400618e79060f9a24fe61fc457aa184b7/e2e/28c8dcfe-f72d-4a1e-8572-60038800e9d0.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/10e063f400618e79060f9a24fe61fc457aa184b7/e2e/ca40efca-d1b0-4c5f-af3c-95fc671a0aee.md"

# ---------------------------------------------------------------------
# Overview sheet: refresh the per-locale status columns (E, F)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.17
$wsOverview.Columns.Item(6).ColumnWidth = 29.17

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("J2").Value = $mdFile1
$wsZh.Range("K2").Value = "28c8dcfe-f72d-4a1e-8572-60038800e9d0.cf4540c688c72ae39d71e35d6c95e3bc092477ac.zh-cn.xlf"
$wsZh.Range("L2").Value = "2016-12-16 09:37:24"

$wsZh.Range("J3").Value = $mdFile2
$wsZh.Range("K3").Value = "ca40efca-d1b0-4c5f-af3c-95fc671a0aee.b9444fcbfd72c09a5391a510e6a8f503855ffd2b.zh-cn.xlf"
$wsZh.Range("L3").Value = "2016-12-16 09:37:24"

$wsZh.Hyperlinks.Add($wsZh.Range("J2"), $mdUrl1, "", "", $mdFile1)
$wsZh.Hyperlinks.Add($wsZh.Range("J3"), $mdUrl2, "", "", $mdFile2)

$wsZh.Range("J2").Font.Underline = 2
$wsZh.Range("J2").Font.Color = 15570276
$wsZh.Range("J3").Font.Underline = 2
$wsZh.Range("J3").Font.Color = 15570276

$wsZh.Columns.Item(3).ColumnWidth = 29.17
$wsZh.Columns.Item(10).ColumnWidth = 39.17
$wsZh.Columns.Item(11).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("J2").Value = $mdFile1
$wsDe.Range("K2").Value = "28c8dcfe-f72d-4a1e-8572-60038800e9d0.cf4540c688c72ae39d71e35d6c95e3bc092477ac.de-de.xlf"
$wsDe.Range("L2").Value = "2016-12-16 09:37:42"

$wsDe.Range("G3").Value = "ca40efca-d1b0-4c5f-af3c-95fc671a0aee.b9444fcbfd72c09a5391a510e6a8f503855ffd2b.de-de.xlf"
$wsDe.Range("J3").Value = $mdFile2
$wsDe.Range("K3").Value = "ca40efca-d1b0-4c5f-af3c-95fc671a0aee.b9444fcbfd72c09a5391a510e6a8f503855ffd2b.de-de.xlf"
$wsDe.Range("L3").Value = "2016-12-16 09:37:42"

$wsDe.Hyperlinks.Add($wsDe.Range("J2"), $mdUrl1, "", "", $mdFile1)
$wsDe.Hyperlinks.Add($wsDe.Range("J3"), $mdUrl2, "", "", $mdFile2)

$wsDe.Range("J2").Font.Underline = 2
$wsDe.Range("J2").Font.Color = 15570276
$wsDe.Range("J3").Font.Underline = 2
$wsDe.Range("J3").Font.Color = 15570276

$wsDe.Columns.Item(3).ColumnWidth = 29.17
$wsDe.Columns.Item(10).ColumnWidth = 39.17
$wsDe.Columns.Item(11).ColumnWidth = 39.17
